# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "64.434.19"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "3.081.13"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.54"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.57"
$ws.Range("E6").Value = "  +5.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.078.04"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.37"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.34"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "3.588.86"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "64.511.18"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "3.088.75"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.77"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.58"
$ws.Range("E23").Value = "  +5.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.58"
$ws.Range("E24").Value = "  +8.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.60"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.25"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.50"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.18"
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.69"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "469.62"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +19.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0836"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0407"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").Value = "2.982.34"
$ws.Range("E41").Value = "  -5.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.30"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.28"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("E46").Value = "  +5.64%  "
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.52"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").Value = "0.0₃0522"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("E51").Value = "  +0.57%  "
